# chore: update Sheets via scheduled runner
# Refreshes market-board derived price/profit columns (H-N) from the
# latest scrape for the affected Leve rows across each crafting sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5025325
$ws.Range("J17").Value = 6699999.5
$ws.Range("L17").Value = 20099998.5
$ws.Range("N17").Value = -20100334.5
$ws.Range("H31").Value = 792.2
$ws.Range("I31").Value = 660.3333
$ws.Range("J31").Value = 990
$ws.Range("K31").Value = 1980.9999
$ws.Range("L31").Value = 2970
$ws.Range("M31").Value = -1750.9999
$ws.Range("N31").Value = -3430
$ws.Range("H111").Value = 1996
$ws.Range("I111").Value = 1993
$ws.Range("J111").Value = 1999
$ws.Range("K111").Value = 5979
$ws.Range("L111").Value = 5997
$ws.Range("M111").Value = -2912
$ws.Range("N111").Value = -12131
$ws.Range("H115").Value = 900.5
$ws.Range("I115").Value = 900.5
$ws.Range("K115").Value = 2701.5
$ws.Range("M115").Value = -1134.5
$ws.Range("H135").Value = 1737.4117
$ws.Range("I135").Value = 1200.0769
$ws.Range("J135").Value = 3483.75
$ws.Range("K135").Value = 10800.6921
$ws.Range("L135").Value = 31353.75
$ws.Range("M135").Value = -8265.6921
$ws.Range("N135").Value = -36423.75
$ws.Range("H137").Value = 2925.853
$ws.Range("I137").Value = 2154.6553
$ws.Range("K137").Value = 6463.965899999999
$ws.Range("M137").Value = -3913.965899999999
$ws.Range("H138").Value = 2752.5686
$ws.Range("I138").Value = 1846.84
$ws.Range("K138").Value = 5540.52
$ws.Range("M138").Value = -400.5199999999995
$ws.Range("H141").Value = 3313.2
$ws.Range("I141").Value = 2014.6666
$ws.Range("K141").Value = 6043.9998
$ws.Range("M141").Value = -863.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11908268
$ws.Range("I32").Value = 7248289.5
$ws.Range("K32").Value = 7248289.5
$ws.Range("M32").Value = -7248002.5
$ws.Range("H132").Value = 2552.0889
$ws.Range("I132").Value = 2180.2058
$ws.Range("J132").Value = 3701.5454
$ws.Range("K132").Value = 6540.617400000001
$ws.Range("L132").Value = 11104.6362
$ws.Range("M132").Value = -4010.617400000001
$ws.Range("N132").Value = -16164.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2419.2
$ws.Range("I105").Value = 1805.16
$ws.Range("K105").Value = 1805.16
$ws.Range("M105").Value = -58.16000000000008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 604.4545000000001
$ws.Range("I7").Value = 564.9
$ws.Range("K7").Value = 564.9
$ws.Range("M7").Value = -451.9
$ws.Range("H8").Value = 8249.75
$ws.Range("J8").Value = 8249.75
$ws.Range("L8").Value = 8249.75
$ws.Range("N8").Value = -8529.75
$ws.Range("H16").Value = 1059.579
$ws.Range("I16").Value = 974.5454999999999
$ws.Range("K16").Value = 974.5454999999999
$ws.Range("M16").Value = -687.5454999999999
$ws.Range("H22").Value = 2625
$ws.Range("I22").Value = 2625
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 2625
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2275
$ws.Range("N22").ClearContents()
$ws.Range("H58").Value = 2450.158
$ws.Range("I58").Value = 1874
$ws.Range("K58").Value = 1874
$ws.Range("M58").Value = -1671
$ws.Range("H86").Value = 4444
$ws.Range("I86").Value = 4444
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4444
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3321
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 4444
$ws.Range("I89").Value = 4444
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 22220
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -16604
$ws.Range("N89").ClearContents()
$ws.Range("H96").Value = 39087.168
$ws.Range("J96").Value = 39087.168
$ws.Range("L96").Value = 39087.168
$ws.Range("N96").Value = -44579.168
$ws.Range("H113").Value = 1059.579
$ws.Range("I113").Value = 974.5454999999999
$ws.Range("K113").Value = 974.5454999999999
$ws.Range("M113").Value = 1195.4545
$ws.Range("H132").Value = 2609.1177
$ws.Range("I132").Value = 2299.4443
$ws.Range("J132").Value = 2957.5
$ws.Range("K132").Value = 6898.3329
$ws.Range("L132").Value = 8872.5
$ws.Range("M132").Value = -4368.3329
$ws.Range("N132").Value = -13932.5
$ws.Range("H136").Value = 2450.158
$ws.Range("I136").Value = 1874
$ws.Range("K136").Value = 5622
$ws.Range("M136").Value = -3072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 899.6667
$ws.Range("I124").Value = 899.6667
$ws.Range("K124").Value = 2699.0001
$ws.Range("M124").Value = 2210.9999
$ws.Range("H140").Value = 13335672
$ws.Range("I140").Value = 13335672
$ws.Range("K140").Value = 40007016
$ws.Range("M140").Value = -40001836

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2556.5715
$ws.Range("I122").Value = 3160.5
$ws.Range("K122").Value = 9481.5
$ws.Range("M122").Value = -7031.5
$ws.Range("H132").Value = 1697.3158
$ws.Range("I132").Value = 1736.3889
$ws.Range("K132").Value = 5209.1667
$ws.Range("M132").Value = -2679.1667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15949.963
$ws.Range("I7").Value = 15620.454
$ws.Range("K7").Value = 15620.454
$ws.Range("M7").Value = -15508.454
$ws.Range("H22").Value = 1516.7826
$ws.Range("J22").Value = 1678.4166
$ws.Range("L22").Value = 1678.4166
$ws.Range("N22").Value = -2268.4166
$ws.Range("H27").Value = 1516.7826
$ws.Range("J27").Value = 1678.4166
$ws.Range("L27").Value = 1678.4166
$ws.Range("N27").Value = -1892.4166
$ws.Range("H40").Value = 8093.1816
$ws.Range("I40").Value = 7902.55
$ws.Range("K40").Value = 7902.55
$ws.Range("M40").Value = -7766.55
$ws.Range("H46").Value = 4273.5
$ws.Range("J46").Value = 4273.5
$ws.Range("L46").Value = 4273.5
$ws.Range("N46").Value = -4649.5
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H75").Value = 110173
$ws.Range("J75").Value = 110173
$ws.Range("L75").Value = 110173
$ws.Range("N75").Value = -112045
$ws.Range("H78").Value = 110173
$ws.Range("J78").Value = 110173
$ws.Range("L78").Value = 330519
$ws.Range("N78").Value = -339879
$ws.Range("H122").Value = 5923.25
$ws.Range("I122").Value = 5923.25
$ws.Range("K122").Value = 17769.75
$ws.Range("M122").Value = -15319.75
$ws.Range("H126").Value = 15949.963
$ws.Range("I126").Value = 15620.454
$ws.Range("K126").Value = 46861.362
$ws.Range("M126").Value = -44391.362
